# Auto-generated Excel COM-interop edit script
# Commit: "July update, incorporate MPI"
# Applies the numeric corrections/updates to existing rows (132-232, 542-630)
# and appends new row 631 (June 2022 data) to Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cellValues = @{
    "H132" = 935363
    "I132" = 100.7
    "M132" = 14066910
    "M133" = 14469363
    "M134" = 14604998
    "M135" = 14690435
    "M136" = 14866894
    "M137" = 14861412
    "M138" = 14670804
    "M139" = 14467029
    "D140" = 823731
    "M140" = 14131809
    "N140" = 8935567
    "M141" = 13894105
    "N141" = 8782082
    "M142" = 13750012
    "N142" = 8776350
    "M143" = 13621881
    "N143" = 8836792
    "N144" = 9042600
    "N145" = 9395220
    "N146" = 9716084
    "N147" = 9916915
    "N148" = 10176548
    "N149" = 10402129
    "N150" = 10586529
    "N151" = 10834442
    "D186" = 590397
    "G186" = 9817199
    "H186" = 1209108
    "M186" = 13889898
    "N186" = 8211202
    "M187" = 13966131
    "N187" = 8231193
    "L188" = 125.89
    "M188" = 14177364
    "N188" = 8310979
    "M189" = 14260672
    "N189" = 8315339
    "M190" = 14422999
    "N190" = 8370920
    "D191" = 573545
    "E191" = 207675
    "G191" = 8670414
    "H191" = 1073055
    "M191" = 14528151
    "N191" = 8368439
    "O191" = 2551090
    "M192" = 14532926
    "N192" = 8290013
    "O192" = 2547712
    "M193" = 14761449
    "N193" = 8365105
    "O193" = 2572360
    "F194" = 14729056
    "M194" = 14844130
    "N194" = 8387992
    "O194" = 2570160
    "M195" = 14889061
    "N195" = 8325464
    "O195" = 2571049
    "F196" = 13097790
    "M196" = 14967742
    "N196" = 8319960
    "O196" = 2563123
    "M197" = 15089887
    "N197" = 8341424
    "O197" = 2566578
    "M198" = 15141503
    "N198" = 8330305
    "O198" = 2567379
    "E199" = 214620
    "F199" = 10187302
    "M199" = 15309047
    "N199" = 8356553
    "O199" = 2587923
    "E200" = 247051
    "M200" = 15440657
    "N200" = 8350572
    "O200" = 2599976
    "M201" = 15510431
    "N201" = 8368315
    "O201" = 2605553
    "F202" = 10071212
    "M202" = 15673741
    "N202" = 8381872
    "O202" = 2635466
    "F203" = 9851479
    "O203" = 2650964
    "O204" = 2658786
    "O205" = 2687335
    "O206" = 2682389
    "O207" = 2694740
    "O208" = 2719323
    "O209" = 2707765
    "O210" = 2680330
    "O211" = 2671978
    "B232" = 1361140
    "F232" = 11589052
    "R542" = 259965
    "R543" = 257775
    "R544" = 261012
    "R545" = 270673
    "R546" = 271406
    "R547" = 277417
    "R548" = 278205
    "R549" = 282646
    "R550" = 276631
    "R551" = 268209
    "R552" = 264976
    "R553" = 275959
    "R554" = 280469
    "R555" = 277366
    "R556" = 277856
    "R557" = 273930
    "R558" = 268982
    "R559" = 277867
    "R560" = 271197
    "R561" = 270605
    "R562" = 275330
    "R563" = 277190
    "R564" = 279346
    "R565" = 273242
    "R566" = 265817
    "R567" = 274813
    "R568" = 274938
    "R569" = 269068
    "R570" = 276032
    "R571" = 276214
    "R572" = 275407
    "R573" = 274823
    "R574" = 276182
    "R575" = 280369
    "R576" = 282255
    "R577" = 280943
    "R578" = 277852
    "R579" = 280562
    "R580" = 283384
    "R581" = 289594
    "R582" = 290459
    "R583" = 292179
    "R584" = 294239
    "R585" = 299237
    "R586" = 288223
    "R587" = 287127
    "R588" = 284038
    "R589" = 285939
    "R590" = 294337
    "R591" = 304478
    "R592" = 306575
    "R593" = 325704
    "R594" = 319802
    "R595" = 317784
    "R596" = 323496
    "R597" = 320574
    "R598" = 326564
    "R599" = 321586
    "R600" = 322969
    "R601" = 327985
    "R602" = 340651
    "R603" = 343283
    "R604" = 346029
    "R605" = 339150
    "R606" = 347971
    "R607" = 342763
    "R608" = 335639
    "R609" = 334827
    "R610" = 335975
    "R611" = 337708
    "R612" = 340627
    "R613" = 343583
    "R614" = 339314
    "R615" = 330248
    "R616" = 331515
    "R617" = 330197
    "R618" = 330699
    "R619" = 319412
    "R620" = 316105
    "R621" = 319627
    "R622" = 313283
    "R623" = 314809
    "R624" = 318834
    "R625" = 316775
    "R626" = 318072
    "R627" = 320065
    "G628" = 6442421
    "H628" = 2542352
    "I628" = 407.11
    "M628" = 35831022
    "R628" = 319759
    "M629" = 33650347
    "P629" = 5218
    "Q629" = 14091
    "R629" = 319650
    "C630" = 1349.6
    "M630" = 31548472
    "P630" = 5254
    "Q630" = 14104
    "R630" = 317251
    "A631" = 44742
    "B631" = 937491
    "D631" = 327971
    "E631" = 106846
    "F631" = 5529454
    "G631" = 4387469
    "H631" = 1694212
    "I631" = 394.7
    "J631" = 31.59
    "K631" = 16.95
    "L631" = 372.74
    "M631" = 29255479
    "N631" = 4860179
    "O631" = 2409837
    "P631" = 5253
    "Q631" = 14109
}

foreach ($cellRef in $cellValues.Keys) {
    $ws.Range($cellRef).Value = $cellValues[$cellRef]
}

Write-Output "Applied $($cellValues.Count) cell updates."
